$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.135.85'
$ws.Range("E2").Value = '  +0.57%  '
$ws.Range("D3").Value = '2.932.72'
$ws.Range("E3").Value = '  +1.02%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.57%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.12'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  +4.62%  '
$ws.Range("E10").Value = '  -0.48%  '
$ws.Range("E11").Value = '  -0.34%  '
$ws.Range("E12").Value = '  -0.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '33.75'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.92%  '
$ws.Range("E14").Value = '  +0.12%  '
$ws.Range("D15").Value = '3.418.27'
$ws.Range("E15").Value = '  +0.89%  '
$ws.Range("D16").Value = '61.055.95'
$ws.Range("E16").Value = '  +0.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.74'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.43%  '
$ws.Range("D18").Value = '2.931.67'
$ws.Range("E18").Value = '  +0.88%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '434.44'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.16%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.680'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.32%  '
$ws.Range("E22").Value = '  +0.38%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '81.60'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("E24").Value = '  +2.43%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.87'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.97%  '
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("E28").Value = '  -1.21%  '
$ws.Range("E29").Value = '  -0.35%  '
$ws.Range("E30").Value = '  -0.88%  '
$ws.Range("E31").Value = '  +2.59%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.76'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.87%  '
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").Value = '0.0₃0870'
$ws.Range("E34").Value = '  +1.76%  '
$ws.Range("E35").Value = '  +0.83%  '
$ws.Range("E36").Value = '  +0.85%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.97'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.21%  '
$ws.Range("E38").Value = '  -0.70%  '
$ws.Range("E39").Value = '  -0.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.61'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '42.16'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.03%  '
$ws.Range("E42").Value = '  -2.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0346'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '373.14'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D45").Value = '2.711.50'
$ws.Range("E45").Value = '  +0.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '133.35'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.39%  '
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.87'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.63%  '
$ws.Range("E49").Value = '  -0.92%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.01'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.12%  '
$ws.Range("E51").Value = '  -0.31%  '
